$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CBBT")

$ws.Range("D17").Value = 4500
$ws.Range("H21").Value = "NA"
$ws.Range("D26").Value = -4900
$ws.Range("D27").Value = -4900
$ws.Range("D33").Value = -4900
$ws.Range("D35").Value = -4900
$ws.Range("D81").Value = -4900
$ws.Range("D91:I91").Value = 0
